# Re-order the weekly price rows (rows 3-20) of the Chirimoya sheet.
# The underlying data for each row (columns D, L, M, N, O, P, Q, S, T -
# the rest are constant across all rows) gets shuffled into a new row,
# per the mapping below (new row -> old row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row -> old row
$mapping = @{
    3  = 5
    4  = 16
    5  = 17
    6  = 4
    7  = 3
    8  = 15
    9  = 8
    10 = 12
    11 = 10
    12 = 11
    13 = 6
    14 = 7
    15 = 13
    16 = 14
    17 = 19
    18 = 9
    19 = 20
    20 = 18
}

# Columns that actually vary row to row.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")

# First, snapshot the current ("before") values for every relevant cell,
# so that writes to one row never clobber data we still need to read
# for another row later.
$snapshot = @{}
foreach ($row in 3..20) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

# Now write the shuffled data back out.
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $rowData = $snapshot[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $rowData[$col]
    }
}
